$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the "blog" entries in row 8 one slot older, and insert a brand new
# blog entry (ser: 120) into the newest slot (C8).
$ws.Range("I8").Value = "type: blog`nwidth: 2`nheight: 1`nser: 118"
$ws.Range("E8").Value = "type: blog`nwidth: 2`nheight: 1`nser: 119"
$ws.Range("C8").Value = "type: blog`nwidth: 2`nheight: 1`nser: 120"
